$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 193 (shifts Dx28's future neighbours - L01.. - down by one)
$ws.Rows("193:193").Insert()

# Populate the new "Dx28 / Hodgkin" row
$ws.Range("A193").Value = "Dx28"
$ws.Range("B193").Value = "Hodgkin"
$ws.Range("C193").Value = "Cancer types"
$ws.Range("D193").Value = "Hodgkin lymphoma"
$ws.Range("E193").Value = "0 = No; 1 = Yes"

# Grow Table1 so the new row becomes part of the table (autoFilter/sort follow along)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E290"))

# Keep the view roughly in sync with where the sheet was scrolled to/selected
$ws.Range("A193").Select()

Write-Output "done"
